$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: fill in previously-empty odds cells with numeric values
$ws.Range("G3").Value = 1.91
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 4.3
$ws.Range("L3").Value = 1.47
$ws.Range("M3").Value = 2.32
$ws.Range("N3").Value = 2.4
$ws.Range("O3").Value = 1.45
$ws.Range("P3").Value = 1.52
$ws.Range("Q3").Value = 2.22
$ws.Range("R3").Value = 2.1
$ws.Range("S3").Value = 1.57
$ws.Range("T3").Value = 5.1
$ws.Range("U3").Value = 7.6
$ws.Range("V3").Value = 9
$ws.Range("W3").Value = 16
$ws.Range("X3").Value = 19.5
$ws.Range("Y3").Value = 45
$ws.Range("Z3").Value = 6.4
$ws.Range("AA3").Value = 6.1
$ws.Range("AB3").Value = 19.5
$ws.Range("AC3").Value = 120
$ws.Range("AE3").Value = 9.25
$ws.Range("AF3").Value = 22
$ws.Range("AG3").Value = 15
$ws.Range("AH3").Value = 75
$ws.Range("AI3").Value = 50
$ws.Range("AJ3").Value = 65

# Row 8: updated odds
$ws.Range("J8").Value = 1.05
$ws.Range("K8").Value = 11
$ws.Range("L8").Value = 1.29
$ws.Range("N8").Value = 1.98
$ws.Range("O8").Value = 1.83
$ws.Range("P8").Value = 1.36

# Row 9: updated odds
$ws.Range("P9").Value = 1.33
$ws.Range("AA9").Value = 7.5

# Row 10: updated odds
$ws.Range("N10").Value = 1.95
$ws.Range("O10").Value = 1.85

# Row 11: updated odds
$ws.Range("G11").Value = 2.6
$ws.Range("H11").Value = 3.4
$ws.Range("I11").Value = 2.63
$ws.Range("J11").Value = 1.05
$ws.Range("K11").Value = 11
$ws.Range("U11").Value = 13
$ws.Range("Z11").Value = 10
$ws.Range("AH11").Value = 26
$ws.Range("AI11").Value = 21

# Row 12: updated odds
$ws.Range("J12").Value = 1.04
$ws.Range("L12").Value = 1.25
$ws.Range("P12").Value = 1.33

# Row 13: updated odds
$ws.Range("G13").Value = 2.2
$ws.Range("H13").Value = 3.4
$ws.Range("L13").Value = 1.25
$ws.Range("M13").Value = 3.75
$ws.Range("N13").Value = 1.88
$ws.Range("O13").Value = 1.93
$ws.Range("AG13").Value = 12
$ws.Range("AI13").Value = 26
$ws.Range("AJ13").Value = 34

# Row 14: updated odds
$ws.Range("G14").Value = 2.05
$ws.Range("I14").Value = 3.3
$ws.Range("P14").Value = 1.33
$ws.Range("R14").Value = 1.62
$ws.Range("T14").Value = 9.5
$ws.Range("AH14").Value = 34
